$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.737.34"
$ws.Range("E2").Value = "  -2.80%  "
$ws.Range("D3").Value = "3.177.54"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.56%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.176.77"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.35%  "
$ws.Range("E10").Value = "  -5.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.476"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.03%  "
$ws.Range("D15").Value = "3.693.68"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").Value = "64.688.87"
$ws.Range("E16").Value = "  -2.93%  "
$ws.Range("D17").Value = "3.174.31"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.717"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.80%  "
$ws.Range("E29").Value = "  -5.32%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.56%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.69%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.66%  "
$ws.Range("E35").Value = "  -5.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("E38").Value = "  +5.80%  "
$ws.Range("D39").Value = "0.0₃0738"
$ws.Range("E39").Value = "  -4.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "457.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0404"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.125"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").Value = "2.868.73"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("E46").Value = "  -7.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.91%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.70%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.116"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.72%  "
